$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2 values (Client Id, User Name, Exam Password, First Name, Last Name)
$ws.Range("A2").Value = "KZtcy695"
$ws.Range("C2").Value = "bdcgwbt47"
$ws.Range("D2").Value = "hJk7!3#R"
$ws.Range("F2").Value = "jjcBWIZR"
$ws.Range("G2").Value = "hukc"

# Update Candidate ID (numeric) in B2
$ws.Range("B2").Value = 23081625
